$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "fixing matricula of Matc65"
# Column A ("matricula") values for rows 14-39 were corrected: every
# row's matricula effectively shifted down by one from where it used to
# be (with row 14 receiving a brand-new matricula, and a few rows in the
# middle of the block being re-ordered), per the target diff.
$values = @{
    14 = "217216526"
    15 = "216117974"
    16 = "221117463"
    17 = "217125254"
    18 = "219218129"
    19 = "218215397"
    20 = "220117282"
    21 = "219217429"
    22 = "216216087"
    23 = "220121412"
    24 = "210201260"
    25 = "201520233"
    26 = "217117994"
    27 = "219118481"
    28 = "221119218"
    29 = "219215012"
    30 = "219121541"
    31 = "214007731"
    32 = "219215013"
    33 = "220117290"
    34 = "219118473"
    35 = "220117273"
    36 = "220120071"
    37 = "221216783"
    38 = "214120645"
    39 = "220217140"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("A$row")
    # The matricula is an ID that must stay text (it was authored as an
    # inline string), not auto-coerced to a number. Forcing the cell to
    # the Text number format while assigning keeps the value a string;
    # resetting back to the Normal style afterwards avoids leaving any
    # stray formatting behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
    $cell.Style = "Normal"
}
